$d = $word.ActiveDocument

# --- Edit 1: expand the opening sentence with a red "(This is a change ... )" note ---
# "This is a Microsoft word document." -> "This is a Microsoft word document.  " (trailing
# double space) followed by three new runs, all colored red, spelling out the parenthetical.
$firstPara = $d.Paragraphs(1).Range
$firstPara.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, `
    $false, $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

$p1 = $d.Paragraphs(1).Range
$ins = $d.Range($p1.End - 1, $p1.End - 1)

$ins.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$ins.Font.Color = 255
$ins.Collapse(0)

$ins.InsertAfter("rsion for main branch")
$ins.Font.Color = 255
$ins.Collapse(0)

$ins.InsertAfter(")")
$ins.Font.Color = 255
$ins.Collapse(0)

# --- Edit 2: drop the trailing "ank God almighty, we are free at last." paragraph ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Delete()

# --- Edit 3: prune the now-unused styles (mirrors Word's own cleanup on save). These
# must be removed from the end of the Styles collection backward, since removing an
# earlier entry first shifts later indices out from under any still-linked styles.
$staleStyles = @(
    "podcast-toolssubscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading4Char",
    "Heading2Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading4",
    "Heading2"
)
foreach ($styleName in $staleStyles) {
    try {
        $d.Styles.Item($styleName).Delete()
    } catch {
        # Style already absent / not removable in this runtime - ignore and continue.
    }
}
